$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the solicitation date for row 8 (Ketal Supermercados)
$ws.Range("E8").Value = 44607

# Update the contact person and email addresses for row 8
$ws.Range("L8").Value = "Jose Luis Murillo"
$ws.Range("M8").Value = "jmurillo@ketal.com.bo; rcondori@ketal.com.bo"

# Give E9 the same (empty) date-styled, underlined formatting as G10
$ws.Range("E9").NumberFormat = "dd/mm/yy;@"
$ws.Range("E9").Font.Underline = $true

# Move the active selection to E9
[void]$ws.Range("E9").Select()
